# "Added school to database"
# Insert a new "school" column (column B) into each of the three grade
# sheets and populate it with the school each team belongs to.

$wb = $excel.ActiveWorkbook

# --- grade78 (sheet2): SMG / NPMG / OMG / SMG -------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Columns.Item(2).Insert()
$ws2.Cells.Item(1, 2).Value = "SMG"
$ws2.Cells.Item(2, 2).Value = "NPMG"
$ws2.Cells.Item(3, 2).Value = "OMG"
$ws2.Cells.Item(4, 2).Value = "SMG"

# --- grade56 (sheet1): SMG / SMG / PMG / PMG / SMG --------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Columns.Item(2).Insert()
$ws1.Cells.Item(1, 2).Value = "SMG"
$ws1.Cells.Item(2, 2).Value = "SMG"
$ws1.Cells.Item(3, 2).Value = "PMG"
$ws1.Cells.Item(4, 2).Value = "PMG"
$ws1.Cells.Item(5, 2).Value = "SMG"

# --- grade912 (sheet3): SMG -------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Columns.Item(2).Insert()
$ws3.Cells.Item(1, 2).Value = "SMG"

# --- restore the selections / active sheet left behind by the edit ---
$ws3.Range("D3").Select()
$ws2.Range("E5").Select()
$ws1.Range("B6").Select()
